$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraped Price column uses "."-grouped / fixed-decimal-point text
# (e.g. "29.379.00", "81.80", "0.00000000120") that Excel would otherwise
# reinterpret as a number and mangle (dropping trailing zeros, switching to
# scientific notation, etc). Force each rewritten Price cell to Text format
# first so the literal string survives untouched, then look/read back as the
# exact scraped text.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.379.00'
$ws.Range("E2").Value = '  +0.18%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.841.75'
$ws.Range("E3").Value = '  -0.10%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  +0.17%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.30'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6258'
$ws.Range("E6").Value = '  -0.23%  '

# Row 7
$ws.Range("E7").Value = '  +0.27%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07426'
$ws.Range("E8").Value = '  -0.88%  '

# Row 9
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2893'
$ws.Range("E9").Value = '  -0.19%  '

# Row 10
$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.93'
$ws.Range("E10").Value = '  +2.15%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07721'
$ws.Range("E11").Value = '  -0.09%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.846.02'
$ws.Range("E12").Value = '  +0.16%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.971'
$ws.Range("E13").Value = '  -0.28%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6742'
$ws.Range("E14").Value = '  -0.79%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001025'
$ws.Range("E15").Value = '  -2.31%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '81.80'
$ws.Range("E16").Value = '  -0.20%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.204'
$ws.Range("E17").Value = '  +0.40%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.382.92'
$ws.Range("E18").Value = '  -0.02%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '233.54'
$ws.Range("E19").Value = '  +2.03%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.31'
$ws.Range("E20").Value = '  -0.04%  '

# Row 21
$ws.Range("E21").Value = '  +0.23%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.271'
$ws.Range("E22").Value = '  -2.92%  '

# Row 23
$ws.Range("E23").Value = '  +0.31%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.32'
$ws.Range("E24").Value = '  -0.10%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.474'
$ws.Range("E25").Value = '  +0.63%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1346'
$ws.Range("E26").Value = '  -1.77%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.30'
$ws.Range("E27").Value = '  -1.16%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.07217'
$ws.Range("E28").Value = '  +12.67%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.463'
$ws.Range("E29").Value = '  +3.99%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.478'
$ws.Range("E30").Value = '  -0.15%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.049'
$ws.Range("E31").Value = '  -0.98%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.023'
$ws.Range("E32").Value = '  -1.64%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.818'
$ws.Range("E33").Value = '  -0.79%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.139'
$ws.Range("E34").Value = '  -0.02%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6971'
$ws.Range("E35").Value = '  +0.04%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.572'
$ws.Range("E36").Value = '  -0.24%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01840'
$ws.Range("E37").Value = '  +0.49%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.914'
$ws.Range("E38").Value = '  +2.89%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.814'
$ws.Range("E39").Value = '  -0.90%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.231.67'
$ws.Range("E40").Value = '  -2.90%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9649'
$ws.Range("E41").Value = '  +5.56%  '

# Row 42
$ws.Range("E42").Value = '  +0.22%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.995.68'
$ws.Range("E43").Value = '  -0.47%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.02'
$ws.Range("E44").Value = '  -0.12%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.43'
$ws.Range("E45").Value = '  -1.07%  '

# Row 46
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000120'
$ws.Range("E46").Value = '  +3.19%  '

# Row 47
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.719'
$ws.Range("E47").Value = '  -0.12%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.947'
$ws.Range("E48").Value = '  -1.86%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.905'
$ws.Range("E49").Value = '  -0.92%  '

# Row 50
$ws.Range("E50").Value = '  -2.23%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3898'
$ws.Range("E51").Value = '  -1.47%  '
